# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
# (Week-17 game stats appended/accumulated into the Seahawks "Team Data" sheets.)

$wb = $excel.ActiveWorkbook

# --- YDS --- (per-play R/P yardage logs: append this week's values)
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value() + " 3 -1 1 20 15 13 23 1 6 4 2 5 5 2 37 7 3 3 6 13 8 0 3 1 2 9 1 6 17 7 0 -2 0 6 3 4 1 30"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 12 7 4 8 9 -5 58 11 28 6 7 13 8 1 13 11 12 17 1 15"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " -1 6 2 1 1 2 0 3 -3 3 26 31 5 0 -1 3 2 0 14 0 1"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 6 12 15 13 16 1 6 20 42 7 16 11 5 6 6 8 9 11 23 4 4"

# --- OFF --- (season-to-date offensive totals updated with Week 17 numbers)
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 388
$ws.Range("E2").Value = 24
$ws.Range("F2").Value = 105
$ws.Range("G2").Value = 111
$ws.Range("J2").Value = 50
$ws.Range("L2").Value = 528
$ws.Range("M2").Value = 357
$ws.Range("O2").Value = 28
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 1031
$ws.Range("C3").Value = 351
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 69
$ws.Range("F3").Value = 187
$ws.Range("G3").Value = 89
$ws.Range("H3").Value = 43
$ws.Range("I3").Value = 104
$ws.Range("J3").Value = 106
$ws.Range("N3").Value = 55

# --- DEF --- (season-to-date defensive totals updated with Week 17 numbers)
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 402
$ws.Range("D2").Value = 29
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = 130
$ws.Range("G2").Value = 130
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 74
$ws.Range("L2").Value = 650
$ws.Range("M2").Value = 418
$ws.Range("O2").Value = 54
$ws.Range("P2").Value = 22
$ws.Range("Q2").Value = 1170
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 452
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 244
$ws.Range("G3").Value = 83
$ws.Range("I3").Value = 130
$ws.Range("J3").Value = 107

# --- ST --- (special teams totals + per-kick logs updated with Week 17)
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 164
$ws.Range("D2").Value = 148
$ws.Range("F2").Value = 260
$ws.Range("G2").Value = 234
$ws.Range("J2").Value = 91
$ws.Range("K2").Value = 85
$ws.Range("L2").Value = 61
$ws.Range("M2").Value = 52
$ws.Range("N2").Value = 37
$ws.Range("O2").Value = 23
$ws.Range("B3").Value = 78
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 57 65"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 47 27"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 21 27"
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 45"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 3 0 30"

# --- TURNS --- (turnovers updated with Week 17)
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 13
$ws.Range("E3").Value = 11

# --- PEN --- (penalties updated with Week 17)
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 35
$ws.Range("D3").Value = 12
